$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $d = $cell.Value2
    $dt = [DateTime]::FromOADate($d)
    $eom = $dt.AddMonths(1).AddDays(-1)
    $cell.Value = $eom.ToOADate()
    $cell.NumberFormat = "m/d/yyyy"
}
